$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "65.50") must be forced
# to stay text-typed (matching the source data, which stores every Price/Volume
# cell as a string) -- otherwise COM auto-converts "65.50" -> 65.5 and drops
# the trailing zero / thousands-dot formatting. Trick: set NumberFormat to Text,
# assign the value, then reset .Style to "Normal" so no stray number-format
# style lingers on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '37.483.08'
$ws.Range('E2').Value = '  +5.42%  '
$ws.Range('D3').Value = '2.054.22'
$ws.Range('E3').Value = '  +3.80%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue 'D5' '252.49'
$ws.Range('E5').Value = '  +2.99%  '
$ws.Range('E6').Value = '  +3.19%  '
Set-TextValue 'D7' '65.50'
$ws.Range('E7').Value = '  +15.13%  '
$ws.Range('E9').Value = '  +7.00%  '
Set-TextValue 'D10' '59.77'
$ws.Range('E10').Value = '  +2.50%  '
Set-TextValue 'D11' '0.0770'
$ws.Range('E11').Value = '  +5.06%  '
$ws.Range('E12').Value = '  +1.61%  '
Set-TextValue 'D13' '0.920'
$ws.Range('E13').Value = '  -2.23%  '
Set-TextValue 'D14' '14.85'
$ws.Range('E14').Value = '  +2.72%  '
Set-TextValue 'D15' '22.64'
$ws.Range('E15').Value = '  +26.23%  '
$ws.Range('D16').Value = '2.355.32'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('E17').Value = '  +6.20%  '
$ws.Range('D18').Value = '2.057.01'
$ws.Range('E18').Value = '  +3.98%  '
$ws.Range('D19').Value = '37.361.31'
$ws.Range('E19').Value = '  +5.26%  '
Set-TextValue 'D20' '73.68'
$ws.Range('E20').Value = '  +3.34%  '
$ws.Range('D21').Value = '0.0₃0877'
$ws.Range('E21').Value = '  +4.37%  '
$ws.Range('E22').Value = '  +6.37%  '
Set-TextValue 'D23' '240.29'
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D24' '2.63'
$ws.Range('E24').Value = '  +5.01%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D25' '1.00'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  +4.84%  '
$ws.Range('E27').Value = '  +11.93%  '
Set-TextValue 'D28' '162.21'
$ws.Range('E28').Value = '  -1.28%  '
Set-TextValue 'D29' '20.02'
$ws.Range('E29').Value = '  +4.79%  '
$ws.Range('E30').Value = '  +23.39%  '
Set-TextValue 'D31' '5.29'
$ws.Range('E31').Value = '  +8.98%  '
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  +9.20%  '
$ws.Range('E34').Value = '  +9.29%  '
$ws.Range('E35').Value = '  +6.34%  '
$ws.Range('E36').Value = '  +1.98%  '
$ws.Range('E37').Value = '  +0.05%  '
Set-TextValue 'D38' '1.83'
$ws.Range('E38').Value = '  +4.28%  '
Set-TextValue 'D39' '6.04'
$ws.Range('E39').Value = '  +16.72%  '
$ws.Range('E40').Value = '  +32.94%  '
$ws.Range('E41').Value = '  +18.10%  '
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range('E43').Value = '  +4.54%  '
$ws.Range('E44').Value = '  +6.19%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D45' '17.26'
$ws.Range('E45').Value = '  +8.27%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D46' '0.0219'
$ws.Range('E46').Value = '  +4.33%  '
Set-TextValue 'D47' '96.80'
$ws.Range('E47').Value = '  +6.05%  '
Set-TextValue 'D48' '7.98'
$ws.Range('E48').Value = '  +5.97%  '
$ws.Range('D49').Value = '1.422.44'
$ws.Range('E49').Value = '  +3.89%  '
$ws.Range('E50').Value = '  +2.34%  '
Set-TextValue 'D51' '46.75'
$ws.Range('E51').Value = '  -1.15%  '
